$wb = $excel.ActiveWorkbook

# ---- sheet: full ----
$ws = $wb.Worksheets.Item("full")
$ws.Range("A2").Value = "Total Rows"
$ws.Range("B2").Value = 999.0
$ws.Range("C2").Value = 1.0
$ws.Range("A3").Value = "Equal Rows"
$ws.Range("B3").Value = 604.0
$ws.Range("C3").Value = 0.6046046046046046
$ws.Range("A4").Value = "Different Rows"
$ws.Range("B4").Value = 395.0
$ws.Range("C4").Value = 0.3953953953953954
$ws.Range("A5").Value = "Matching Rows"
$ws.Range("B5").Value = 960.0
$ws.Range("C5").Value = 0.960960960960961
$ws.Range("A6").Value = "Non-matching Rows"
$ws.Range("B6").Value = 960.0
$ws.Range("C6").Value = 0.03903903903903904
$ws.Range("A7").Value = "(Source1) - Total Rows"
$ws.Range("B7").Value = 967.0
$ws.Range("C7").Value = 0.9679679679679679
$ws.Range("A8").Value = "(Source1) - Matching Rows"
$ws.Range("B8").Value = 7.0
$ws.Range("C8").Value = 0.9927611168562565
$ws.Range("A9").Value = "(Source1) - Non-matching Rows"
$ws.Range("B9").Value = 7.0
$ws.Range("C9").Value = 0.007238883143743537
$ws.Range("A10").Value = "(Source2) - Total Rows"
$ws.Range("B10").Value = 992.0
$ws.Range("C10").Value = 0.992992992992993
$ws.Range("A11").Value = "(Source2) - Matching Rows"
$ws.Range("B11").Value = 32.0
$ws.Range("C11").Value = 0.967741935483871
$ws.Range("A12").Value = "(Source2) - Non-matching Rows"
$ws.Range("B12").Value = 32.0
$ws.Range("C12").Value = 0.03225806451612903

# ---- sheet: left ----
$ws = $wb.Worksheets.Item("left")
$ws.Range("A2").Value = "Total Rows"
$ws.Range("B2").Value = 967.0
$ws.Range("C2").Value = 1.0
$ws.Range("A3").Value = "Equal Rows"
$ws.Range("B3").Value = 604.0
$ws.Range("C3").Value = 0.6246122026887281
$ws.Range("A4").Value = "Different Rows"
$ws.Range("B4").Value = 363.0
$ws.Range("C4").Value = 0.375387797311272
$ws.Range("A5").Value = "Matching Rows"
$ws.Range("B5").Value = 960.0
$ws.Range("C5").Value = 0.9927611168562565
$ws.Range("A6").Value = "Non-matching Rows"
$ws.Range("B6").Value = 960.0
$ws.Range("C6").Value = 0.007238883143743537
$ws.Range("A7").Value = "(Source1) - Total Rows"
$ws.Range("B7").Value = 967.0
$ws.Range("C7").Value = 1.0
$ws.Range("A8").Value = "(Source1) - Matching Rows"
$ws.Range("B8").Value = 7.0
$ws.Range("C8").Value = 0.9927611168562565
$ws.Range("A9").Value = "(Source1) - Non-matching Rows"
$ws.Range("B9").Value = 7.0
$ws.Range("C9").Value = 0.007238883143743537
$ws.Range("A10").Value = "(Source2) - Total Rows"
$ws.Range("B10").Value = 960.0
$ws.Range("C10").Value = 0.9927611168562565
$ws.Range("A11").Value = "(Source2) - Matching Rows"
$ws.Range("B11").Value = 0.0
$ws.Range("C11").Value = 1.0
$ws.Range("A12").Value = "(Source2) - Non-matching Rows"
$ws.Range("B12").Value = 0.0
$ws.Range("C12").Value = 0.0

# ---- sheet: right ----
$ws = $wb.Worksheets.Item("right")
$ws.Range("A2").Value = "Total Rows"
$ws.Range("B2").Value = 992.0
$ws.Range("C2").Value = 1.0
$ws.Range("A3").Value = "Equal Rows"
$ws.Range("B3").Value = 604.0
$ws.Range("C3").Value = 0.6088709677419355
$ws.Range("A4").Value = "Different Rows"
$ws.Range("B4").Value = 388.0
$ws.Range("C4").Value = 0.3911290322580645
$ws.Range("A5").Value = "Matching Rows"
$ws.Range("B5").Value = 960.0
$ws.Range("C5").Value = 0.967741935483871
$ws.Range("A6").Value = "Non-matching Rows"
$ws.Range("B6").Value = 960.0
$ws.Range("C6").Value = 0.03225806451612903
$ws.Range("A7").Value = "(Source1) - Total Rows"
$ws.Range("B7").Value = 960.0
$ws.Range("C7").Value = 0.967741935483871
$ws.Range("A8").Value = "(Source1) - Matching Rows"
$ws.Range("B8").Value = 0.0
$ws.Range("C8").Value = 1.0
$ws.Range("A9").Value = "(Source1) - Non-matching Rows"
$ws.Range("B9").Value = 0.0
$ws.Range("C9").Value = 0.0
$ws.Range("A10").Value = "(Source2) - Total Rows"
$ws.Range("B10").Value = 992.0
$ws.Range("C10").Value = 1.0
$ws.Range("A11").Value = "(Source2) - Matching Rows"
$ws.Range("B11").Value = 32.0
$ws.Range("C11").Value = 0.967741935483871
$ws.Range("A12").Value = "(Source2) - Non-matching Rows"
$ws.Range("B12").Value = 32.0
$ws.Range("C12").Value = 0.03225806451612903

# ---- sheet: inner ----
$ws = $wb.Worksheets.Item("inner")
$ws.Range("A2").Value = "Total Rows"
$ws.Range("B2").Value = 960.0
$ws.Range("C2").Value = 1.0
$ws.Range("A3").Value = "Equal Rows"
$ws.Range("B3").Value = 604.0
$ws.Range("C3").Value = 0.6291666666666667
$ws.Range("A4").Value = "Different Rows"
$ws.Range("B4").Value = 356.0
$ws.Range("C4").Value = 0.37083333333333335
$ws.Range("A5").Value = "Matching Rows"
$ws.Range("B5").Value = 960.0
$ws.Range("C5").Value = 1.0
$ws.Range("A6").Value = "Non-matching Rows"
$ws.Range("B6").Value = 960.0
$ws.Range("C6").Value = 0.0
$ws.Range("A7").Value = "(Source1) - Total Rows"
$ws.Range("B7").Value = 960.0
$ws.Range("C7").Value = 1.0
$ws.Range("A8").Value = "(Source1) - Matching Rows"
$ws.Range("B8").Value = 0.0
$ws.Range("C8").Value = 1.0
$ws.Range("A9").Value = "(Source1) - Non-matching Rows"
$ws.Range("B9").Value = 0.0
$ws.Range("C9").Value = 0.0
$ws.Range("A10").Value = "(Source2) - Total Rows"
$ws.Range("B10").Value = 960.0
$ws.Range("C10").Value = 1.0
$ws.Range("A11").Value = "(Source2) - Matching Rows"
$ws.Range("B11").Value = 0.0
$ws.Range("C11").Value = 1.0
$ws.Range("A12").Value = "(Source2) - Non-matching Rows"
$ws.Range("B12").Value = 0.0
$ws.Range("C12").Value = 0.0

# ---- sheet: complement ----
$ws = $wb.Worksheets.Item("complement")
$ws.Range("A2").Value = "Total Rows"
$ws.Range("B2").Value = 39.0
$ws.Range("C2").Value = 1.0
$ws.Range("A3").Value = "Equal Rows"
$ws.Range("B3").Value = 0.0
$ws.Range("C3").Value = 0.0
$ws.Range("A4").Value = "Different Rows"
$ws.Range("B4").Value = 39.0
$ws.Range("C4").Value = 1.0
$ws.Range("A5").Value = "Matching Rows"
$ws.Range("B5").Value = 0.0
$ws.Range("C5").Value = 0.0
$ws.Range("A6").Value = "Non-matching Rows"
$ws.Range("B6").Value = 0.0
$ws.Range("C6").Value = 1.0
$ws.Range("A7").Value = "(Source1) - Total Rows"
$ws.Range("B7").Value = 7.0
$ws.Range("C7").Value = 0.1794871794871795
$ws.Range("A8").Value = "(Source1) - Matching Rows"
$ws.Range("B8").Value = 7.0
$ws.Range("C8").Value = 0.0
$ws.Range("A9").Value = "(Source1) - Non-matching Rows"
$ws.Range("B9").Value = 7.0
$ws.Range("C9").Value = 1.0
$ws.Range("A10").Value = "(Source2) - Total Rows"
$ws.Range("B10").Value = 32.0
$ws.Range("C10").Value = 0.8205128205128205
$ws.Range("A11").Value = "(Source2) - Matching Rows"
$ws.Range("B11").Value = 32.0
$ws.Range("C11").Value = 0.0
$ws.Range("A12").Value = "(Source2) - Non-matching Rows"
$ws.Range("B12").Value = 32.0
$ws.Range("C12").Value = 1.0

# ---- sheet: lcomp ----
$ws = $wb.Worksheets.Item("lcomp")
$ws.Range("A2").Value = "Total Rows"
$ws.Range("B2").Value = 7.0
$ws.Range("C2").Value = 1.0
$ws.Range("A3").Value = "Equal Rows"
$ws.Range("B3").Value = 0.0
$ws.Range("C3").Value = 0.0
$ws.Range("A4").Value = "Different Rows"
$ws.Range("B4").Value = 7.0
$ws.Range("C4").Value = 1.0
$ws.Range("A5").Value = "Matching Rows"
$ws.Range("B5").Value = 0.0
$ws.Range("C5").Value = 0.0
$ws.Range("A6").Value = "Non-matching Rows"
$ws.Range("B6").Value = 0.0
$ws.Range("C6").Value = 1.0
$ws.Range("A7").Value = "(Source1) - Total Rows"
$ws.Range("B7").Value = 7.0
$ws.Range("C7").Value = 1.0
$ws.Range("A8").Value = "(Source1) - Matching Rows"
$ws.Range("B8").Value = 7.0
$ws.Range("C8").Value = 0.0
$ws.Range("A9").Value = "(Source1) - Non-matching Rows"
$ws.Range("B9").Value = 7.0
$ws.Range("C9").Value = 1.0
$ws.Range("A10").Value = "(Source2) - Total Rows"
$ws.Range("B10").Value = 0.0
$ws.Range("C10").Value = 0.0
$ws.Range("A11").Value = "(Source2) - Matching Rows"
$ws.Range("B11").Value = 0.0
$ws.Range("C11").Value = "<NaN>"
$ws.Range("A12").Value = "(Source2) - Non-matching Rows"
$ws.Range("B12").Value = 0.0
$ws.Range("C12").Value = "<NaN>"

# ---- sheet: rcomp ----
$ws = $wb.Worksheets.Item("rcomp")
$ws.Range("A2").Value = "Total Rows"
$ws.Range("B2").Value = 32.0
$ws.Range("C2").Value = 1.0
$ws.Range("A3").Value = "Equal Rows"
$ws.Range("B3").Value = 0.0
$ws.Range("C3").Value = 0.0
$ws.Range("A4").Value = "Different Rows"
$ws.Range("B4").Value = 32.0
$ws.Range("C4").Value = 1.0
$ws.Range("A5").Value = "Matching Rows"
$ws.Range("B5").Value = 0.0
$ws.Range("C5").Value = 0.0
$ws.Range("A6").Value = "Non-matching Rows"
$ws.Range("B6").Value = 0.0
$ws.Range("C6").Value = 1.0
$ws.Range("A7").Value = "(Source1) - Total Rows"
$ws.Range("B7").Value = 0.0
$ws.Range("C7").Value = 0.0
$ws.Range("A8").Value = "(Source1) - Matching Rows"
$ws.Range("B8").Value = 0.0
$ws.Range("C8").Value = "<NaN>"
$ws.Range("A9").Value = "(Source1) - Non-matching Rows"
$ws.Range("B9").Value = 0.0
$ws.Range("C9").Value = "<NaN>"
$ws.Range("A10").Value = "(Source2) - Total Rows"
$ws.Range("B10").Value = 32.0
$ws.Range("C10").Value = 1.0
$ws.Range("A11").Value = "(Source2) - Matching Rows"
$ws.Range("B11").Value = 32.0
$ws.Range("C11").Value = 0.0
$ws.Range("A12").Value = "(Source2) - Non-matching Rows"
$ws.Range("B12").Value = 32.0
$ws.Range("C12").Value = 1.0

# ---- update AutoFilter ranges ----
$ws = $wb.Worksheets.Item("full")
$ws.AutoFilterMode = $false
$ws.Range("A1:C12").AutoFilter()
$ws = $wb.Worksheets.Item("left")
$ws.AutoFilterMode = $false
$ws.Range("A1:C12").AutoFilter()
$ws = $wb.Worksheets.Item("right")
$ws.AutoFilterMode = $false
$ws.Range("A1:C12").AutoFilter()
$ws = $wb.Worksheets.Item("inner")
$ws.AutoFilterMode = $false
$ws.Range("A1:C12").AutoFilter()
$ws = $wb.Worksheets.Item("complement")
$ws.AutoFilterMode = $false
$ws.Range("A1:C12").AutoFilter()
$ws = $wb.Worksheets.Item("lcomp")
$ws.AutoFilterMode = $false
$ws.Range("A1:C12").AutoFilter()
$ws = $wb.Worksheets.Item("rcomp")
$ws.AutoFilterMode = $false
$ws.Range("A1:C12").AutoFilter()

# ---- update _FilterDatabase defined names ----
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
  $nm = $names.Item($i)
  if ($nm.Name -like "*_FilterDatabase") {
    $sheetPart = $nm.Name.Split("!")[0]
    $nm.RefersTo = "=" + $sheetPart + "!`$A`$1:`$C`$12"
  }
}
